{"js": "// Align the financial-statement template with the budget guidelines:\n//  - move the \"_GoBack\" bookmark from the very first (empty) paragraph to\n//    the end of the \"By Ekhagastiftelsen granted sum\" row\n//  - relabel / reorder the rows of the Expenses table\n\n// --- 1. Move the _GoBack bookmark ------------------------------------------\n// Originally it sits in the first, empty paragraph of the document.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Table 1 (index 0) = Application number / Project title / ...\n// Table 2 (index 1) = By Ekhagastiftelsen granted sum\n// Table 3 (index 2) = Expenses\n// Table 4 (index 3) = Place and date / Signature\nconst sumTable = tables.items[1];\nconst expensesTable = tables.items[2];\n\n// Re-create the bookmark at the end of the last cell of the \"granted sum\"\n// table, i.e. right after its visible content.\nconst lastCell = sumTable.getCell(0, 1);\nconst goBackRange = lastCell.body.getRange(\"End\");\ngoBackRange.insertBookmark(\"_GoBack\");\n\n// --- 2. Update the Expenses table ------------------------------------------\nexpensesTable.getCell(1, 0).value = \"Salaries\";\nexpensesTable.getCell(2, 0).value = \"Equipment, instruments, materials\";\nexpensesTable.getCell(3, 0).value = \"Travel\";\nexpensesTable.getCell(4, 0).value = \"Office and building costs\";\nexpensesTable.getCell(5, 0).value = \"Other\";\n\nawait context.sync();\n", "ps1": "# Align the financial-statement template with the budget guidelines:\n#  - move the \"_GoBack\" bookmark from the very first (empty) paragraph to\n#    the end of the \"By Ekhagastiftelsen granted sum\" row\n#  - relabel / reorder the rows of the Expenses table\n\n$d = $word.ActiveDocument\n\n# --- 1. Move the _GoBack bookmark -----------------------------------------\n# Originally it sits in the first, empty paragraph of the document.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Re-create it at the end of the last cell of the \"granted sum\" table\n# (Table 2), i.e. right after its visible content.\n$sumTable = $d.Tables(2)\n$lastCell = $sumTable.Rows(1).Cells($sumTable.Rows(1).Cells.Count)\n$goBackRange = $lastCell.Range\n$goBackRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n\n# --- 2. Update the Expenses table ------------------------------------------\n$expensesTable = $d.Tables(3)\n\n$expensesTable.Rows(2).Cells(1).Range.Text = \"Salaries\"\n$expensesTable.Rows(3).Cells(1).Range.Text = \"Equipment, instruments, materials\"\n$expensesTable.Rows(4).Cells(1).Range.Text = \"Travel\"\n$expensesTable.Rows(5).Cells(1).Range.Text = \"Office and building costs\"\n$expensesTable.Rows(6).Cells(1).Range.Text = \"Other\"\n"}
